$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BLNK")

# Row 14 - Non Recurring
$ws.Range("D14").Value = 20900
$ws.Range("E14").Value = 0

# Row 17 - Total Operating Expenses
$ws.Range("D17").Value = 30500
$ws.Range("E17").Value = 10500
$ws.Range("F17").Value = 15300
$ws.Range("G17").Value = 29500

# Row 18 - Operating Income or Loss
$ws.Range("D18").Value = -28000
$ws.Range("E18").Value = -7200

# Row 20 - Total Other Income/Expenses Net
$ws.Range("D20").Value = -44100
$ws.Range("E20").Value = 700

# Row 32 - Other Items
$ws.Range("D32").Value = 44100
$ws.Range("E32").Value = -700
